# Weekly update: a new price record is added for "Femacal de La Calera - Papaya".
# Insert a new row above the current row 32, pushing the existing rows 32:50
# down to 33:51, then populate the new row 32 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(32).Insert()

$ws.Cells.Item(32, 1).Value  = 3
$ws.Cells.Item(32, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(32, 3).Value  = "Coquimbo"
$ws.Cells.Item(32, 4).Value  = 44603
$ws.Cells.Item(32, 5).Value  = 5
$ws.Cells.Item(32, 6).Value  = "Fruta"
$ws.Cells.Item(32, 7).Value  = 100108
$ws.Cells.Item(32, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(32, 9).Value  = 100108004
$ws.Cells.Item(32, 10).Value = "Papaya"
$ws.Cells.Item(32, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 48
$ws.Cells.Item(32, 14).Value = 23000
$ws.Cells.Item(32, 15).Value = 23000
$ws.Cells.Item(32, 16).Value = 23000
$ws.Cells.Item(32, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(32, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(32, 19).Value = 2300
$ws.Cells.Item(32, 20).Value = 10
